# "Bis Kapitel 5.2.2 gekommen"
# Log two more days of work on the "Schriftliche Arbeit" / "Fertigstellung
# der Arbeit" line (rows 121-122 of the Stundenerfassung sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 121 already exists (only A121 had a date-formatted but empty cell) -
# fill in the rest of the row.
$ws.Cells.Item(121, 1).Value = 42971
$ws.Cells.Item(121, 2).Value = "Schriftliche Arbeit"
$ws.Cells.Item(121, 3).Value = "Fertigstellung der Arbeit"
$ws.Cells.Item(121, 4).Value = 7

# Row 122 is brand new. Clone the date number format from A121 (so the new
# cell keeps the existing "m/d/yyyy"-style formatting instead of picking up
# the bare column default) and then fill in the values.
$ws.Range("A121").Copy()
$ws.Range("A122").PasteSpecial(-4122)
$ws.Cells.Item(122, 1).Value = 42972
$ws.Cells.Item(122, 2).Value = "Schriftliche Arbeit"
$ws.Cells.Item(122, 3).Value = "Fertigstellung der Arbeit"
$ws.Cells.Item(122, 4).Value = 6

# Scroll the visible window down a bit so the new rows stay in view
# (topLeftCell moved from A100 to A105 in the saved view state).
$win = $excel.ActiveWindow
$win.ScrollRow = 105
$win.ScrollColumn = 1
